# edit.ps1 - applies the "submission 3" revision to First Meeting.docx
#
# Two independent changes, per the commit's unified diff:
#   1. Drop the (now stale) spell-check markers <w:proofErr .../> that
#      bracket the "GearShare" run in the Project Name row.
#   2. Re-split the single "Maheen Siddique, Linton Dsouza, Cyrus Chakma"
#      run into three runs (".../C" + "i" + "rus Chakma") with identical
#      run formatting - same visible text, just authored as separate runs
#      (e.g. from retyping the "i" in "Cyrus").
#
# Both paragraphs are located at runtime via Find (rather than a
# hard-coded paragraph index) and then rewritten in place with
# Range.InsertXML using a minimal WordProcessingML package fragment that
# reproduces the original paragraph/run properties exactly, so nothing
# besides the targeted markup changes.

$d = $word.ActiveDocument

function Set-ParagraphXml($paragraph, [string]$innerBodyXml) {
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' + $innerBodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $paragraph.Range.InsertXML($pkg)
}

# --- 1. Remove proofErr spellStart/spellEnd around "GearShare" ----------

$find1 = $d.Content
$find1.Find.Execute("GearShare", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
if (-not $find1.Find.Found) {
    throw "Could not locate the 'GearShare' paragraph"
}
$gearSharePara = $find1.Paragraphs(1)

$gearShareXml = '<w:p w14:paraId="19AA0A62" w14:textId="19C964A2" w:rsidR="00954B3E" w:rsidRPr="006A59CE" w:rsidRDefault="006F4AB9" w:rsidP="0034601D">' +
    '<w:pPr>' +
        '<w:widowControl w:val="0"/>' +
        '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
            '<w:color w:val="000000"/>' +
            '<w:lang w:eastAsia="en-CA"/>' +
        '</w:rPr>' +
    '</w:pPr>' +
    '<w:r>' +
        '<w:rPr>' +
            '<w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
            '<w:color w:val="000000"/>' +
            '<w:lang w:eastAsia="en-CA"/>' +
        '</w:rPr>' +
        '<w:t>GearShare</w:t>' +
    '</w:r>' +
'</w:p>'

Set-ParagraphXml $gearSharePara $gearShareXml

# --- 2. Re-split the "Cyrus" run in the invitees list --------------------

$find2 = $d.Content
$find2.Find.Execute("Cyrus Chakma", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
if (-not $find2.Find.Found) {
    throw "Could not locate the 'Cyrus Chakma' paragraph"
}
$inviteesPara = $find2.Paragraphs(1)

$runPr = '<w:rPr>' +
            '<w:rFonts w:ascii="Arial" w:eastAsia="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
            '<w:color w:val="000000"/>' +
            '<w:lang w:eastAsia="en-CA"/>' +
         '</w:rPr>'

$inviteesXml = '<w:p w14:paraId="58783887" w14:textId="77777777" w:rsidR="00954B3E" w:rsidRDefault="00954B3E" w:rsidP="0034601D">' +
    '<w:pPr>' +
        '<w:widowControl w:val="0"/>' +
        '<w:pBdr>' +
            '<w:top w:val="nil"/>' +
            '<w:left w:val="nil"/>' +
            '<w:bottom w:val="nil"/>' +
            '<w:right w:val="nil"/>' +
            '<w:between w:val="nil"/>' +
        '</w:pBdr>' +
        '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>' +
        $runPr +
    '</w:pPr>' +
    '<w:r>' + $runPr + '<w:t>Maheen Siddique, Linton Dsouza, C</w:t></w:r>' +
    '<w:r>' + $runPr + '<w:t>i</w:t></w:r>' +
    '<w:r>' + $runPr + '<w:t>rus Chakma</w:t></w:r>' +
'</w:p>'

Set-ParagraphXml $inviteesPara $inviteesXml

Write-Host "Done."
